# Changed actual energy file
# - Convert A43/D43 from text to real numbers
# - Insert a duplicate row 44 (same values as the corrected row 43)
# - Insert a new row 45 with text-typed values (A45="27", D45="20")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A reference cell that carries the workbook's default (unstyled) cell style,
# used later to strip any incidental number-formatting Excel may apply when a
# date-looking string is typed into a cell.
$defaultStyle = $ws.Range("A2").Style

# --- Fix row 43: A43 and D43 should be real numbers, not text ---
$ws.Range("A43").Value = 27
$ws.Range("D43").Value = 19

# --- Row 44: duplicate of the corrected row 43 ---
$ws.Range("A44").Value = 27
$ws.Range("B44").Value = "Partly Cloudy"

$c44 = $ws.Cells.Item(44, 3)
$c44.NumberFormat = "@"
$c44.Value = "01/18/2025"
$c44.Style = $defaultStyle

$ws.Range("D44").Value = 19

# --- Row 45: new row, values stored as text (matches source data) ---
$a45 = $ws.Cells.Item(45, 1)
$a45.NumberFormat = "@"
$a45.Value = "27"
$a45.Style = $defaultStyle

$ws.Range("B45").Value = "Partly Cloudy"

$c45 = $ws.Cells.Item(45, 3)
$c45.NumberFormat = "@"
$c45.Value = "01/18/2025"
$c45.Style = $defaultStyle

$d45 = $ws.Cells.Item(45, 4)
$d45.NumberFormat = "@"
$d45.Value = "20"
$d45.Style = $defaultStyle
